$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 181 ---
# Date (col A) - copy the date/time style from the row above so it reuses
# the existing numFmt/style instead of minting a new one.
$ws.Cells.Item(180, 1).Copy() | Out-Null
$ws.Cells.Item(181, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(181, 1).Value = 45457.2916666667

$ws.Cells.Item(181, 2).Value = 0
$ws.Cells.Item(181, 3).Value = 5.15000009536743
$ws.Cells.Item(181, 4).Value = 5.15000009536743
$ws.Cells.Item(181, 5).Value = 5.15000009536743
$ws.Cells.Item(181, 6).Value = 5.15000009536743

# adj_close (col G) must stay text, so force a text format before writing it,
# then restore the default (Normal) style so no stray "@" format sticks
# around on the cell itself.
$ws.Cells.Item(181, 7).NumberFormat = "@"
$ws.Cells.Item(181, 7).Value = "5.15000009536743"
$ws.Cells.Item(181, 7).Style = "Normal"

$ws.Cells.Item(181, 8).Value = "VLC.MI"

# --- Row 182 ---
$ws.Cells.Item(180, 1).Copy() | Out-Null
$ws.Cells.Item(182, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(182, 1).Value = 45460.4668981482

$ws.Cells.Item(182, 2).Value = 1080
$ws.Cells.Item(182, 3).Value = 5.05000019073486
$ws.Cells.Item(182, 4).Value = 4.98000001907349
$ws.Cells.Item(182, 5).Value = 4.98000001907349
$ws.Cells.Item(182, 6).Value = 5.05000019073486

$ws.Cells.Item(182, 7).NumberFormat = "@"
$ws.Cells.Item(182, 7).Value = "5.05000019073486"
$ws.Cells.Item(182, 7).Style = "Normal"

$ws.Cells.Item(182, 8).Value = "VLC.MI"

$excel.CutCopyMode = 0
